$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "63.801.34"
Set-TextCell "E2" "  +0.38%  "
Set-TextCell "D3" "2.626.39"
Set-TextCell "E3" "  -0.81%  "
Set-TextCell "E4" "  -0.07%  "
Set-TextCell "D5" "596.15"
Set-TextCell "E5" "  -0.86%  "
Set-TextCell "D6" "150.05"
Set-TextCell "E6" "  +2.11%  "
Set-TextCell "E7" "  -0.08%  "
Set-TextCell "E8" "  +0.01%  "
Set-TextCell "E9" "  +0.58%  "
Set-TextCell "D10" "5.69"
Set-TextCell "E10" "  +1.85%  "
Set-TextCell "D11" "0.382"
Set-TextCell "E11" "  +3.44%  "
Set-TextCell "E12" "  -1.18%  "
Set-TextCell "D13" "27.72"
Set-TextCell "E13" "  +0.94%  "
Set-TextCell "D14" "3.095.99"
Set-TextCell "E14" "  -1.00%  "
Set-TextCell "D15" "63.577.42"
Set-TextCell "E15" "  +0.17%  "
Set-TextCell "E16" "  +2.99%  "
Set-TextCell "D17" "2.661.65"
Set-TextCell "E17" "  +0.03%  "
Set-TextCell "D18" "12.31"
Set-TextCell "E18" "  +7.50%  "
Set-TextCell "D19" "4.64"
Set-TextCell "E19" "  +2.15%  "
Set-TextCell "D20" "348.71"
Set-TextCell "E21" "  -1.24%  "
Set-TextCell "D22" "0.999"
Set-TextCell "E22" "  -0.17%  "
Set-TextCell "D23" "5.69"
Set-TextCell "E23" "  +1.98%  "
Set-TextCell "D24" "66.32"
Set-TextCell "E24" "  -0.59%  "
Set-TextCell "D25" "1.74"
Set-TextCell "E25" "  +12.89%  "
Set-TextCell "D26" "9.20"
Set-TextCell "E26" "  +1.56%  "
Set-TextCell "E27" "  -1.03%  "
Set-TextCell "D28" "562.70"
Set-TextCell "E28" "  +0.41%  "
Set-TextCell "D29" "8.24"
Set-TextCell "E29" "  +4.87%  "
Set-TextCell "E30" "  +0.55%  "
Set-TextCell "E31" "  +0.19%  "
Set-TextCell "D32" "2.04"
Set-TextCell "E32" "  +1.47%  "
Set-TextCell "D33" "0.0₃0845"
Set-TextCell "E33" "  +3.73%  "
Set-TextCell "E34" "  -0.21%  "
Set-TextCell "D35" "5.21"
Set-TextCell "E35" "  +1.10%  "
Set-TextCell "D36" "168.85"
Set-TextCell "E36" "  +0.96%  "
Set-TextCell "E37" "  +0.58%  "
Set-TextCell "D38" "1.00"
Set-TextCell "E38" "  -0.16%  "
Set-TextCell "E39" "  -0.04%  "
Set-TextCell "D40" "19.34"
Set-TextCell "E40" "  +1.33%  "
Set-TextCell "B41" "USDe"
Set-TextCell "C41" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D41" "0.998"
Set-TextCell "E41" "  -0.04%  "
Set-TextCell "B42" "Aave"
Set-TextCell "C42" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D42" "169.85"
Set-TextCell "E42" "  +0.93%  "
Set-TextCell "D43" "39.89"
Set-TextCell "E43" "  -0.13%  "
Set-TextCell "D44" "3.91"
Set-TextCell "E44" "  +3.95%  "
Set-TextCell "D45" "0.0595"
Set-TextCell "E45" "  +4.08%  "
Set-TextCell "D46" "21.34"
Set-TextCell "E46" "  -3.60%  "
Set-TextCell "E47" "  +0.07%  "
Set-TextCell "D48" "0.0248"
Set-TextCell "E48" "  +0.50%  "
Set-TextCell "D49" "1.99"
Set-TextCell "E49" "  +6.65%  "
Set-TextCell "D50" "0.0967"
Set-TextCell "E50" "  +0.79%  "
Set-TextCell "E51" "  +2.19%  "
